$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E for these rows are stored as text (inlineStr) in the source file.
# We force a Text number format before assigning so Excel COM does not coerce
# numeric-looking strings (e.g. "8.09") into real numbers, then reset the style
# back to Normal so no stray formatting is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "60.278.00"
Set-TextValue $ws.Range("E2") "  -0.40%  "
Set-TextValue $ws.Range("D3") "2.385.36"
Set-TextValue $ws.Range("E3") "  -0.96%  "
Set-TextValue $ws.Range("E4") "  +0.43%  "
Set-TextValue $ws.Range("D5") "563.17"
Set-TextValue $ws.Range("E5") "  -1.14%  "
Set-TextValue $ws.Range("D6") "139.62"
Set-TextValue $ws.Range("E6") "  +0.69%  "
Set-TextValue $ws.Range("E7") "  -0.25%  "
Set-TextValue $ws.Range("D8") "0.534"
Set-TextValue $ws.Range("E8") "  +1.78%  "
Set-TextValue $ws.Range("D9") "2.386.75"
Set-TextValue $ws.Range("E9") "  -0.11%  "
Set-TextValue $ws.Range("E10") "  -1.73%  "
Set-TextValue $ws.Range("D11") "0.158"
Set-TextValue $ws.Range("E11") "  -0.73%  "
Set-TextValue $ws.Range("D12") "5.14"
Set-TextValue $ws.Range("E12") "  +1.64%  "
Set-TextValue $ws.Range("D13") "0.341"
Set-TextValue $ws.Range("E13") "  +0.85%  "
Set-TextValue $ws.Range("D14") "25.90"
Set-TextValue $ws.Range("E14") "  +0.59%  "
Set-TextValue $ws.Range("D15") "2.834.69"
Set-TextValue $ws.Range("E15") "  -0.53%  "
Set-TextValue $ws.Range("E16") "  -1.71%  "
Set-TextValue $ws.Range("D17") "60.196.46"
Set-TextValue $ws.Range("E17") "  -0.51%  "
Set-TextValue $ws.Range("D18") "2.395.21"
Set-TextValue $ws.Range("E18") "  +0.09%  "
Set-TextValue $ws.Range("D19") "8.09"
Set-TextValue $ws.Range("E19") "  +11.35%  "
Set-TextValue $ws.Range("D20") "10.54"
Set-TextValue $ws.Range("E20") "  -0.27%  "
Set-TextValue $ws.Range("D21") "322.89"
Set-TextValue $ws.Range("E21") "  +1.00%  "
Set-TextValue $ws.Range("D22") "4.08"
Set-TextValue $ws.Range("E22") "  +1.62%  "
Set-TextValue $ws.Range("D23") "6.01"
Set-TextValue $ws.Range("E23") "  -0.68%  "
Set-TextValue $ws.Range("E24") "  -0.13%  "
Set-TextValue $ws.Range("E25") "  -2.36%  "
Set-TextValue $ws.Range("D26") "64.51"
Set-TextValue $ws.Range("E26") "  +0.13%  "
Set-TextValue $ws.Range("D27") "564.37"
Set-TextValue $ws.Range("E27") "  -1.16%  "
Set-TextValue $ws.Range("D28") "8.10"
Set-TextValue $ws.Range("E28") "  -3.81%  "
Set-TextValue $ws.Range("D29") "2.501.38"
Set-TextValue $ws.Range("E29") "  +0.24%  "
Set-TextValue $ws.Range("D30") "0.0₃0931"
Set-TextValue $ws.Range("E30") "  +1.45%  "
Set-TextValue $ws.Range("D31") "8.03"
Set-TextValue $ws.Range("E31") "  +2.30%  "
Set-TextValue $ws.Range("E32") "  -1.52%  "
Set-TextValue $ws.Range("D33") "1.81"
Set-TextValue $ws.Range("E33") "  -1.32%  "
Set-TextValue $ws.Range("E34") "  -0.99%  "
Set-TextValue $ws.Range("E35") "  -0.59%  "
Set-TextValue $ws.Range("E36") "  +5.41%  "
Set-TextValue $ws.Range("D37") "153.07"
Set-TextValue $ws.Range("E37") "  +2.65%  "
Set-TextValue $ws.Range("E38") "  -0.07%  "
Set-TextValue $ws.Range("D39") "0.369"
Set-TextValue $ws.Range("E39") "  +0.24%  "
Set-TextValue $ws.Range("D40") "18.23"
Set-TextValue $ws.Range("E40") "  +0.33%  "
Set-TextValue $ws.Range("D41") "5.11"
Set-TextValue $ws.Range("E41") "  +0.19%  "
Set-TextValue $ws.Range("E42") "  -0.08%  "
Set-TextValue $ws.Range("D43") "41.67"
Set-TextValue $ws.Range("E43") "  +1.78%  "
Set-TextValue $ws.Range("D45") "2.48"
Set-TextValue $ws.Range("E45") "  +5.52%  "
Set-TextValue $ws.Range("D46") "0.0₆0279"
Set-TextValue $ws.Range("E46") "  -2.89%  "
Set-TextValue $ws.Range("D47") "141.21"
Set-TextValue $ws.Range("E47") "  +0.65%  "
Set-TextValue $ws.Range("D48") "3.54"
Set-TextValue $ws.Range("E48") "  +1.05%  "
Set-TextValue $ws.Range("D49") "0.589"
Set-TextValue $ws.Range("E49") "  +0.69%  "
Set-TextValue $ws.Range("D50") "0.0504"
Set-TextValue $ws.Range("E50") "  +0.37%  "
Set-TextValue $ws.Range("D51") "19.23"
Set-TextValue $ws.Range("E51") "  +0.08%  "

Write-Host "Updated cryptos list"
